$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update match results for the two teams involved in the postponed ("aplazado")
# jornada 9 games. Row 3 (team previously "Coquina") and Row 5 (team previously
# "Armada") both get new stats, and their team names are swapped.
$ws.Range("A3").Value = "Armada"
$ws.Range("B3").Value = 18
$ws.Range("D3").Value = 49
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 0

$ws.Range("A5").Value = "Coquina"
$ws.Range("B5").Value = 13
$ws.Range("D5").Value = 67
$ws.Range("E5").Value = 4
$ws.Range("G5").Value = 1

# Re-sort the standings table (A2:H7) by Puntos (column B) descending.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B7"), 0, 2)
$ws.Sort.SetRange($ws.Range("A2:H7"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Move the active selection to F4, matching the recorded UI state.
$ws.Range("F4").Select() | Out-Null
